$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs target cluster)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.020676
$ws.Range("H2").Value = 0.062028
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2679253333333333
$ws.Range("N2").Value = 0.8037759999999999
$ws.Range("O2").Value = 0.1226600350746756
$ws.Range("P2").Value = 0.1226600350746756
$ws.Range("Q2").Value = 0.005539624191999999
$ws.Range("R2").Value = 0.04985661772799999
$ws.Range("S2").Value = 0.1226600350746756
$ws.Range("T2").Value = 0.1226600350746756

# Row 3 (FAPs target cluster)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.020676
$ws.Range("H3").Value = 0.062028
$ws.Range("O3").Value = 0.327101565785771
$ws.Range("P3").Value = 0.327101565785771
$ws.Range("Q3").Value = 0.014772698752
$ws.Range("R3").Value = 0.132954288768
$ws.Range("S3").Value = 0.327101565785771
$ws.Range("T3").Value = 0.327101565785771

# Row 4 (MuSCs target cluster)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.020676
$ws.Range("H4").Value = 0.062028
$ws.Range("M4").Value = 1.145196333333333
$ws.Range("N4").Value = 3.435589
$ws.Range("O4").Value = 0.5242871984759059
$ws.Range("P4").Value = 0.5242871984759059
$ws.Range("Q4").Value = 0.023678079388
$ws.Range("R4").Value = 0.213102714492
$ws.Range("S4").Value = 0.5242871984759059
$ws.Range("T4").Value = 0.5242871984759059

# Row 5 (Resolving-Mac target cluster)
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.020676
$ws.Range("H5").Value = 0.062028
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05668500000000001
$ws.Range("N5").Value = 0.170055
$ws.Range("O5").Value = 0.02595120066364754
$ws.Range("P5").Value = 0.02595120066364754
$ws.Range("Q5").Value = 0.00117201906
$ws.Range("R5").Value = 0.01054817154
$ws.Range("S5").Value = 0.02595120066364754
$ws.Range("T5").Value = 0.02595120066364754
